$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new TODO item just before "Mark physical axes with blue and
#    red." The new item is added as a sibling ListParagraph (same numbering
#    list/level) right after "Make README.md for the firmware.".
# ---------------------------------------------------------------------------
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ppText = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ppText -eq "Make README.md for the firmware.") {
        $anchorIndex = $i
        break
    }
}

$anchor = $d.Paragraphs.Item($anchorIndex)
$anchor.Range.InsertParagraphAfter()

$newParaIndex = $anchorIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)

$firstPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$firstPoint.InsertAfter("Uninstall + reinstall Ubuntu and make a")

$newParaRange = $d.Paragraphs.Item($newParaIndex).Range
$secondPoint = $d.Range($newParaRange.End - 1, $newParaRange.End - 1)
$secondPoint.InsertAfter("n install_requirements.txt file by making a list of necessary packages")

# ---------------------------------------------------------------------------
# 2. Mark three already-finished TODO items as done by striking them through.
# ---------------------------------------------------------------------------
$strikeTargets = @(
    "Make separate GitHub for foam machine.",
    "Update this GitHub repository with the latest information.",
    "Make separate, public, GitHub."
)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $text = $pp.Range.Text.TrimEnd([char]13, [char]7)
    if ($strikeTargets -contains $text) {
        $pp.Range.Font.StrikeThrough = 1
    }
}
